$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh. D-column cells are kept as
# plain text (matching the source inline-string cells) by forcing a
# text number format before the write, so Excel does not silently
# coerce values like "0.590" or "1.00" into numbers and drop the
# trailing zeroes.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.299.15"
$ws.Range("E2").Value = "  +5.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.464.55"
$ws.Range("E3").Value = "  +6.48%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "566.43"
$ws.Range("E5").Value = "  +3.89%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.23"
$ws.Range("E6").Value = "  +10.55%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +2.60%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.462.36"
$ws.Range("E9").Value = "  +6.48%  "

$ws.Range("E10").Value = "  +4.84%  "

$ws.Range("E11").Value = "  +3.16%  "

$ws.Range("E12").Value = "  +1.30%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  +5.42%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.42"
$ws.Range("E14").Value = "  +12.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.909.85"
$ws.Range("E15").Value = "  +6.66%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.245.06"
$ws.Range("E16").Value = "  +4.95%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000142"
$ws.Range("E17").Value = "  +6.45%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.469.01"
$ws.Range("E18").Value = "  +6.79%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.24"
$ws.Range("E19").Value = "  +6.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "341.49"
$ws.Range("E20").Value = "  +9.14%  "

$ws.Range("E21").Value = "  +5.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.79"
$ws.Range("E22").Value = "  +3.51%  "

$ws.Range("E23").Value = "  -0.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.64"
$ws.Range("E24").Value = "  +2.55%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.174"
$ws.Range("E25").Value = "  +2.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"

$ws.Range("E27").Value = "  +9.37%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.11"
$ws.Range("E28").Value = "  +3.62%  "

$ws.Range("E29").Value = "  +8.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.87"
$ws.Range("E30").Value = "  +15.75%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0₃0810"
$ws.Range("E31").Value = "  +11.85%  "

$ws.Range("E32").Value = "  +7.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "175.56"
$ws.Range("E33").Value = "  +2.37%  "

$ws.Range("E34").Value = "  +11.40%  "

$ws.Range("E35").Value = "  +4.78%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.89"
$ws.Range("E36").Value = "  +5.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "370.96"
$ws.Range("E37").Value = "  +17.29%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  +8.66%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.71"
$ws.Range("E41").Value = "  +12.75%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.44"
$ws.Range("E42").Value = "  +6.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "151.43"
$ws.Range("E43").Value = "  +10.87%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.71"
$ws.Range("E44").Value = "  +6.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.58"
$ws.Range("E45").Value = "  +8.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.598"
$ws.Range("E46").Value = "  +6.26%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0961"
$ws.Range("E47").Value = "  +2.49%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  +4.77%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0237"
$ws.Range("E49").Value = "  +10.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0225"
$ws.Range("E50").Value = "  +4.79%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.03"
$ws.Range("E51").Value = "  +7.13%  "

